$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2208.028465307524
$ws.Range("C2").Value = 1298.338748465416
$ws.Range("D2").Value = 952.8587701285444
$ws.Range("E2").Value = 830.5100773001739

$ws.Range("B3").Value = 2051.863742562773
$ws.Range("C3").Value = 1183.083383140179
$ws.Range("D3").Value = 894.6238311929591
$ws.Range("E3").Value = 795.2769952160717

$ws.Range("B4").Value = 2218.78299239892
$ws.Range("C4").Value = 1319.607786861474
$ws.Range("D4").Value = 995.0177761066858
$ws.Range("E4").Value = 890.9557820942919

$ws.Range("B5").Value = 795.3625414706451
$ws.Range("C5").Value = 475.4013437469265
$ws.Range("D5").Value = 390.6279785442163
$ws.Range("E5").Value = 370.479860731409

$ws.Range("B6").Value = 204.0546052485936
$ws.Range("C6").Value = 125.1759797846952
$ws.Range("D6").Value = 90.69209456185537
$ws.Range("E6").Value = 76.02898833840067

$ws.Range("B7").Value = 28.64197684667349
$ws.Range("C7").Value = 18.29898074934333
$ws.Range("D7").Value = 15.00777732763962
$ws.Range("E7").Value = 13.6517290310976

$ws.Range("B8").Value = 2602.825699863444
$ws.Range("C8").Value = 1617.062296026409
$ws.Range("D8").Value = 1151.048916303432
$ws.Range("E8").Value = 945.8714456117903

$ws.Range("B9").Value = 497.3364501225034
$ws.Range("C9").Value = 308.280478403597
$ws.Range("D9").Value = 212.7743839416862
$ws.Range("E9").Value = 172.5767032627227

$ws.Range("B10").Value = 389.6230887972482
$ws.Range("C10").Value = 262.5000129926964
$ws.Range("D10").Value = 253.9793968266577
$ws.Range("E10").Value = 255.7124706025495

$ws.Range("B11").Value = 57.99316334920526
$ws.Range("C11").Value = 36.93277136971744
$ws.Range("D11").Value = 36.98338595019786
$ws.Range("E11").Value = 41.33569351649567

$ws.Range("B12").Value = 31.17083543065311
$ws.Range("C12").Value = 24.4611878464834
$ws.Range("D12").Value = 37.26329659219444
$ws.Range("E12").Value = 46.59883366587204

$ws.Range("B13").Value = 151.1694474754235
$ws.Range("C13").Value = 104.0276315578221
$ws.Range("D13").Value = 112.777791526568
$ws.Range("E13").Value = 121.143424118054
